$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# New task rows appended at the bottom of the log (rows 55-57)
$ws.Range("A55").Value = 45359
$ws.Range("B55").Value = 2
$ws.Range("C55").Value = "Bugs fixen"

$ws.Range("A56").Value = 45364
$ws.Range("B56").Value = 5
$ws.Range("C56").Value = "Proposals Review"

$ws.Range("A57").Value = 45366
$ws.Range("B57").Value = 2
$ws.Range("C57").Value = "Ui fixes"

# Match date number formatting used by the rest of column A
$ws.Range("A55:A57").NumberFormat = $ws.Range("A54").NumberFormat

# Update the view so the new rows are visible, mirroring the author's scroll position
$ws.Range("G54").Select()
